# Scheduled-runner price/profit refresh: updates the currentAveragePrice*,
# LevePrice*/LeveProfit* columns (H:N) on a handful of rows across each
# job sheet with freshly scraped market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 42811.43
$ws.Range("I32").Value = 8043.8335
$ws.Range("J32").Value = 68887.125
$ws.Range("K32").Value = 8043.8335
$ws.Range("L32").Value = 68887.125
$ws.Range("M32").Value = -7717.8335
$ws.Range("N32").Value = -69539.125
$ws.Range("H69").Value = 6668.8335
$ws.Range("I69").Value = 5003.25
$ws.Range("K69").Value = 15009.75
$ws.Range("M69").Value = -14135.75
$ws.Range("H70").Value = 921
$ws.Range("I70").Value = 741
$ws.Range("J70").Value = 957
$ws.Range("K70").Value = 2223
$ws.Range("L70").Value = 2871
$ws.Range("M70").Value = -1953
$ws.Range("N70").Value = -3411
$ws.Range("H72").Value = 6668.8335
$ws.Range("I72").Value = 5003.25
$ws.Range("K72").Value = 45029.25
$ws.Range("M72").Value = -40661.25
$ws.Range("H73").Value = 921
$ws.Range("I73").Value = 741
$ws.Range("J73").Value = 957
$ws.Range("K73").Value = 2223
$ws.Range("L73").Value = 2871
$ws.Range("M73").Value = -1287
$ws.Range("N73").Value = -4743
$ws.Range("H80").Value = 4807.4614
$ws.Range("I80").Value = 2249.75
$ws.Range("J80").Value = 5944.222
$ws.Range("K80").Value = 6749.25
$ws.Range("L80").Value = 17832.666
$ws.Range("M80").Value = -5751.25
$ws.Range("N80").Value = -19828.666
$ws.Range("H83").Value = 4807.4614
$ws.Range("I83").Value = 2249.75
$ws.Range("J83").Value = 5944.222
$ws.Range("K83").Value = 20247.75
$ws.Range("L83").Value = 53497.998
$ws.Range("M83").Value = -15255.75
$ws.Range("N83").Value = -63481.998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7539.975
$ws.Range("I32").Value = 4144.4116
$ws.Range("K32").Value = 4144.4116
$ws.Range("M32").Value = -3857.4116
$ws.Range("H49").Value = 52530
$ws.Range("J49").Value = 52530
$ws.Range("L49").Value = 52530
$ws.Range("N49").Value = -53050
$ws.Range("H63").Value = 2127.0908
$ws.Range("I63").Value = 2271.1428
$ws.Range("J63").Value = 1875
$ws.Range("K63").Value = 2271.1428
$ws.Range("L63").Value = 1875
$ws.Range("M63").Value = -1585.1428
$ws.Range("N63").Value = -3247
$ws.Range("H66").Value = 2127.0908
$ws.Range("I66").Value = 2271.1428
$ws.Range("J66").Value = 1875
$ws.Range("K66").Value = 11355.714
$ws.Range("L66").Value = 9375
$ws.Range("M66").Value = -7923.714
$ws.Range("N66").Value = -16239
$ws.Range("H95").Value = 39736
$ws.Range("J95").Value = 39736
$ws.Range("L95").Value = 39736
$ws.Range("N95").Value = -45228
$ws.Range("H97").Value = 1620.1666
$ws.Range("I97").Value = 996.44446
$ws.Range("J97").Value = 3491.3333
$ws.Range("K97").Value = 996.44446
$ws.Range("L97").Value = 3491.3333
$ws.Range("M97").Value = -500.44446
$ws.Range("N97").Value = -4483.3333
$ws.Range("H101").Value = 27801
$ws.Range("J101").Value = 27801
$ws.Range("L101").Value = 27801
$ws.Range("N101").Value = -34291
$ws.Range("H105").Value = 118999.5
$ws.Range("J105").Value = 118999.5
$ws.Range("L105").Value = 118999.5
$ws.Range("N105").Value = -125987.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1993.8096
$ws.Range("I86").Value = 1775.8667
$ws.Range("J86").Value = 2538.6667
$ws.Range("K86").Value = 1775.8667
$ws.Range("L86").Value = 2538.6667
$ws.Range("M86").Value = -652.8667
$ws.Range("N86").Value = -4784.6667
$ws.Range("H89").Value = 1993.8096
$ws.Range("I89").Value = 1775.8667
$ws.Range("J89").Value = 2538.6667
$ws.Range("K89").Value = 8879.333500000001
$ws.Range("L89").Value = 12693.3335
$ws.Range("M89").Value = -3263.333500000001
$ws.Range("N89").Value = -23925.3335
$ws.Range("H94").Value = 1098.2858
$ws.Range("I94").Value = 1037.3334
$ws.Range("K94").Value = 1037.3334
$ws.Range("M94").Value = -586.3334
$ws.Range("H134").Value = 2853.625
$ws.Range("I134").Value = 2923.1
$ws.Range("J134").Value = 2737.8333
$ws.Range("K134").Value = 8769.299999999999
$ws.Range("L134").Value = 8213.499899999999
$ws.Range("M134").Value = -6234.299999999999
$ws.Range("N134").Value = -13283.4999
$ws.Range("H140").Value = 77456.336
$ws.Range("J140").Value = 77456.336
$ws.Range("L140").Value = 77456.336
$ws.Range("N140").Value = -87816.336
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10401.568
$ws.Range("I31").Value = 1828.5518
$ws.Range("J31").Value = 26976.066
$ws.Range("K31").Value = 1828.5518
$ws.Range("L31").Value = 26976.066
$ws.Range("M31").Value = -1533.5518
$ws.Range("N31").Value = -27566.066
$ws.Range("H34").Value = 10401.568
$ws.Range("I34").Value = 1828.5518
$ws.Range("J34").Value = 26976.066
$ws.Range("K34").Value = 1828.5518
$ws.Range("L34").Value = 26976.066
$ws.Range("M34").Value = -1626.5518
$ws.Range("N34").Value = -27380.066
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 1586.4
$ws.Range("I132").Value = 1269.619
$ws.Range("K132").Value = 3808.857
$ws.Range("M132").Value = -1278.857
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 20.375
$ws.Range("I12").Value = 19.5
$ws.Range("K12").Value = 58.5
$ws.Range("M12").Value = 114.5
$ws.Range("H113").Value = 1086.1333
$ws.Range("I113").Value = 958.3333
$ws.Range("J113").Value = 1118.0834
$ws.Range("K113").Value = 2874.9999
$ws.Range("L113").Value = 3354.2502
$ws.Range("M113").Value = -704.9998999999998
$ws.Range("N113").Value = -7694.2502
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 71431470
$ws.Range("I80").Value = 111113850
$ws.Range("K80").Value = 111113850
$ws.Range("M80").Value = -111112852
$ws.Range("H83").Value = 71431470
$ws.Range("I83").Value = 111113850
$ws.Range("K83").Value = 555569250
$ws.Range("M83").Value = -555564258
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4072.4614
$ws.Range("I68").Value = 4327
$ws.Range("K68").Value = 4327
$ws.Range("M68").Value = -3578
$ws.Range("H71").Value = 4072.4614
$ws.Range("I71").Value = 4327
$ws.Range("K71").Value = 21635
$ws.Range("M71").Value = -17891
$ws.Range("H93").Value = 3041.1738
$ws.Range("I93").Value = 2945.4
$ws.Range("K93").Value = 2945.4
$ws.Range("M93").Value = -1697.4
$ws.Range("H106").Value = 19997.5
$ws.Range("J106").Value = 19997.5
$ws.Range("L106").Value = 19997.5
$ws.Range("N106").Value = -22521.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11083
$ws.Range("J62").Value = 11083
$ws.Range("L62").Value = 11083
$ws.Range("N62").Value = -12331
$ws.Range("H65").Value = 11083
$ws.Range("J65").Value = 11083
$ws.Range("L65").Value = 55415
$ws.Range("N65").Value = -61655
$ws.Range("H81").Value = 100
$ws.Range("I81").Value = 100
$ws.Range("K81").Value = 200
$ws.Range("M81").Value = 861
$ws.Range("H84").Value = 100
$ws.Range("I84").Value = 100
$ws.Range("K84").Value = 1000
$ws.Range("M84").Value = 4304
$ws.Range("H94").Value = 46499.3
$ws.Range("I94").Value = 34997.5
$ws.Range("J94").Value = 49374.75
$ws.Range("K94").Value = 34997.5
$ws.Range("L94").Value = 49374.75
$ws.Range("M94").Value = -34096.5
$ws.Range("N94").Value = -51176.75
$ws.Range("H95").Value = 106948.5
$ws.Range("J95").Value = 106948.5
$ws.Range("L95").Value = 106948.5
$ws.Range("N95").Value = -112440.5

Write-Host "Applied all changes"